$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
$ws.Range("D2").Value = "25.986.21"
$ws.Range("D3").Value = "1.641.22"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.15"
$ws.Range("D5").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0638"
$ws.Range("D9").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.870.36"
$ws.Range("D14").Value = "1.643.08"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.546"
$ws.Range("D15").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.90"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "25.942.20"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.99"
$ws.Range("D20").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.53"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.94"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("D29").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0499"
$ws.Range("D31").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.905"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").Value = "1.134.46"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.34"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").Value = "1.779.78"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.69"
$ws.Range("D46").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.73"
$ws.Range("D49").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").ClearFormats()

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("E24").Value = "  +7.64%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  +0.81%  "
